# Append new sensor-log rows to the mmWave sheets (auto update Excel log)
$wb = $excel.ActiveWorkbook

# --- mmWave(InBed): append rows 14-26 ---
$ws = $wb.Worksheets.Item("mmWave(InBed)")
$ws.Cells.Item(14, 1).Value = "'2026-02-01"
$ws.Cells.Item(14, 2).Value = "20:58:42"
$ws.Cells.Item(14, 3).Value = "20:00"
$ws.Cells.Item(14, 4).Value = "Bedroom"
$ws.Cells.Item(14, 5).Value = "In Bed"
$ws.Cells.Item(14, 6).Value = "Occupied"
$ws.Cells.Item(15, 1).Value = "'2026-02-01"
$ws.Cells.Item(15, 2).Value = "20:58:48"
$ws.Cells.Item(15, 3).Value = "20:00"
$ws.Cells.Item(15, 4).Value = "Bedroom"
$ws.Cells.Item(15, 5).Value = "In Bed"
$ws.Cells.Item(15, 6).Value = "Occupied"
$ws.Cells.Item(16, 1).Value = "'2026-02-01"
$ws.Cells.Item(16, 2).Value = "20:58:56"
$ws.Cells.Item(16, 3).Value = "20:00"
$ws.Cells.Item(16, 4).Value = "Bedroom"
$ws.Cells.Item(16, 5).Value = "In Bed"
$ws.Cells.Item(16, 6).Value = "Occupied"
$ws.Cells.Item(17, 1).Value = "'2026-02-01"
$ws.Cells.Item(17, 2).Value = "20:59:01"
$ws.Cells.Item(17, 3).Value = "20:00"
$ws.Cells.Item(17, 4).Value = "Bedroom"
$ws.Cells.Item(17, 5).Value = "In Bed"
$ws.Cells.Item(17, 6).Value = "Occupied"
$ws.Cells.Item(18, 1).Value = "'2026-02-01"
$ws.Cells.Item(18, 2).Value = "20:59:02"
$ws.Cells.Item(18, 3).Value = "20:00"
$ws.Cells.Item(18, 4).Value = "Bedroom"
$ws.Cells.Item(18, 5).Value = "In Bed"
$ws.Cells.Item(18, 6).Value = "Occupied"
$ws.Cells.Item(19, 1).Value = "'2026-02-01"
$ws.Cells.Item(19, 2).Value = "20:59:03"
$ws.Cells.Item(19, 3).Value = "20:00"
$ws.Cells.Item(19, 4).Value = "Bedroom"
$ws.Cells.Item(19, 5).Value = "In Bed"
$ws.Cells.Item(19, 6).Value = "Occupied"
$ws.Cells.Item(20, 1).Value = "'2026-02-01"
$ws.Cells.Item(20, 2).Value = "20:59:04"
$ws.Cells.Item(20, 3).Value = "20:00"
$ws.Cells.Item(20, 4).Value = "Bedroom"
$ws.Cells.Item(20, 5).Value = "In Bed"
$ws.Cells.Item(20, 6).Value = "Occupied"
$ws.Cells.Item(21, 1).Value = "'2026-02-01"
$ws.Cells.Item(21, 2).Value = "20:59:05"
$ws.Cells.Item(21, 3).Value = "20:00"
$ws.Cells.Item(21, 4).Value = "Bedroom"
$ws.Cells.Item(21, 5).Value = "In Bed"
$ws.Cells.Item(21, 6).Value = "Occupied"
$ws.Cells.Item(22, 1).Value = "'2026-02-01"
$ws.Cells.Item(22, 2).Value = "20:59:11"
$ws.Cells.Item(22, 3).Value = "20:00"
$ws.Cells.Item(22, 4).Value = "Bedroom"
$ws.Cells.Item(22, 5).Value = "In Bed"
$ws.Cells.Item(22, 6).Value = "Occupied"
$ws.Cells.Item(23, 1).Value = "'2026-02-01"
$ws.Cells.Item(23, 2).Value = "20:59:16"
$ws.Cells.Item(23, 3).Value = "20:00"
$ws.Cells.Item(23, 4).Value = "Bedroom"
$ws.Cells.Item(23, 5).Value = "In Bed"
$ws.Cells.Item(23, 6).Value = "Occupied"
$ws.Cells.Item(24, 1).Value = "'2026-02-01"
$ws.Cells.Item(24, 2).Value = "20:59:23"
$ws.Cells.Item(24, 3).Value = "20:00"
$ws.Cells.Item(24, 4).Value = "Bedroom"
$ws.Cells.Item(24, 5).Value = "In Bed"
$ws.Cells.Item(24, 6).Value = "Occupied"
$ws.Cells.Item(25, 1).Value = "'2026-02-01"
$ws.Cells.Item(25, 2).Value = "20:59:24"
$ws.Cells.Item(25, 3).Value = "20:00"
$ws.Cells.Item(25, 4).Value = "Bedroom"
$ws.Cells.Item(25, 5).Value = "In Bed"
$ws.Cells.Item(25, 6).Value = "Occupied"
$ws.Cells.Item(26, 1).Value = "'2026-02-01"
$ws.Cells.Item(26, 2).Value = "20:59:25"
$ws.Cells.Item(26, 3).Value = "20:00"
$ws.Cells.Item(26, 4).Value = "Bedroom"
$ws.Cells.Item(26, 5).Value = "In Bed"
$ws.Cells.Item(26, 6).Value = "Occupied"
$ws.Range("A14:A26").ClearFormats()

# --- mmWave(BR): append rows 13-25 ---
$ws = $wb.Worksheets.Item("mmWave(BR)")
$ws.Cells.Item(13, 1).Value = "'2026-02-01"
$ws.Cells.Item(13, 2).Value = "20:58:42"
$ws.Cells.Item(13, 3).Value = "20:00"
$ws.Cells.Item(13, 4).Value = "Bedroom"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = "Occupied"
$ws.Cells.Item(14, 1).Value = "'2026-02-01"
$ws.Cells.Item(14, 2).Value = "20:58:48"
$ws.Cells.Item(14, 3).Value = "20:00"
$ws.Cells.Item(14, 4).Value = "Bedroom"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = "Occupied"
$ws.Cells.Item(15, 1).Value = "'2026-02-01"
$ws.Cells.Item(15, 2).Value = "20:58:56"
$ws.Cells.Item(15, 3).Value = "20:00"
$ws.Cells.Item(15, 4).Value = "Bedroom"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = "Occupied"
$ws.Cells.Item(16, 1).Value = "'2026-02-01"
$ws.Cells.Item(16, 2).Value = "20:59:02"
$ws.Cells.Item(16, 3).Value = "20:00"
$ws.Cells.Item(16, 4).Value = "Bedroom"
$ws.Cells.Item(16, 5).Value = 10
$ws.Cells.Item(16, 6).Value = "Occupied"
$ws.Cells.Item(17, 1).Value = "'2026-02-01"
$ws.Cells.Item(17, 2).Value = "20:59:02"
$ws.Cells.Item(17, 3).Value = "20:00"
$ws.Cells.Item(17, 4).Value = "Bedroom"
$ws.Cells.Item(17, 5).Value = 4
$ws.Cells.Item(17, 6).Value = "Occupied"
$ws.Cells.Item(18, 1).Value = "'2026-02-01"
$ws.Cells.Item(18, 2).Value = "20:59:03"
$ws.Cells.Item(18, 3).Value = "20:00"
$ws.Cells.Item(18, 4).Value = "Bedroom"
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = "Occupied"
$ws.Cells.Item(19, 1).Value = "'2026-02-01"
$ws.Cells.Item(19, 2).Value = "20:59:05"
$ws.Cells.Item(19, 3).Value = "20:00"
$ws.Cells.Item(19, 4).Value = "Bedroom"
$ws.Cells.Item(19, 5).Value = 20
$ws.Cells.Item(19, 6).Value = "Occupied"
$ws.Cells.Item(20, 1).Value = "'2026-02-01"
$ws.Cells.Item(20, 2).Value = "20:59:05"
$ws.Cells.Item(20, 3).Value = "20:00"
$ws.Cells.Item(20, 4).Value = "Bedroom"
$ws.Cells.Item(20, 5).Value = 2
$ws.Cells.Item(20, 6).Value = "Occupied"
$ws.Cells.Item(21, 1).Value = "'2026-02-01"
$ws.Cells.Item(21, 2).Value = "20:59:11"
$ws.Cells.Item(21, 3).Value = "20:00"
$ws.Cells.Item(21, 4).Value = "Bedroom"
$ws.Cells.Item(21, 5).Value = 1
$ws.Cells.Item(21, 6).Value = "Occupied"
$ws.Cells.Item(22, 1).Value = "'2026-02-01"
$ws.Cells.Item(22, 2).Value = "20:59:17"
$ws.Cells.Item(22, 3).Value = "20:00"
$ws.Cells.Item(22, 4).Value = "Bedroom"
$ws.Cells.Item(22, 5).Value = 2
$ws.Cells.Item(22, 6).Value = "Occupied"
$ws.Cells.Item(23, 1).Value = "'2026-02-01"
$ws.Cells.Item(23, 2).Value = "20:59:23"
$ws.Cells.Item(23, 3).Value = "20:00"
$ws.Cells.Item(23, 4).Value = "Bedroom"
$ws.Cells.Item(23, 5).Value = 72
$ws.Cells.Item(23, 6).Value = "Occupied"
$ws.Cells.Item(24, 1).Value = "'2026-02-01"
$ws.Cells.Item(24, 2).Value = "20:59:24"
$ws.Cells.Item(24, 3).Value = "20:00"
$ws.Cells.Item(24, 4).Value = "Bedroom"
$ws.Cells.Item(24, 5).Value = 91
$ws.Cells.Item(24, 6).Value = "Occupied"
$ws.Cells.Item(25, 1).Value = "'2026-02-01"
$ws.Cells.Item(25, 2).Value = "20:59:26"
$ws.Cells.Item(25, 3).Value = "20:00"
$ws.Cells.Item(25, 4).Value = "Bedroom"
$ws.Cells.Item(25, 5).Value = 2
$ws.Cells.Item(25, 6).Value = "Occupied"
$ws.Range("A13:A25").ClearFormats()

# --- mmWave(HR): append rows 13-25 ---
$ws = $wb.Worksheets.Item("mmWave(HR)")
$ws.Cells.Item(13, 1).Value = "'2026-02-01"
$ws.Cells.Item(13, 2).Value = "20:58:42"
$ws.Cells.Item(13, 3).Value = "20:00"
$ws.Cells.Item(13, 4).Value = "Bedroom"
$ws.Cells.Item(13, 5).Value = 50
$ws.Cells.Item(13, 6).Value = "Occupied"
$ws.Cells.Item(14, 1).Value = "'2026-02-01"
$ws.Cells.Item(14, 2).Value = "20:58:48"
$ws.Cells.Item(14, 3).Value = "20:00"
$ws.Cells.Item(14, 4).Value = "Bedroom"
$ws.Cells.Item(14, 5).Value = 49
$ws.Cells.Item(14, 6).Value = "Occupied"
$ws.Cells.Item(15, 1).Value = "'2026-02-01"
$ws.Cells.Item(15, 2).Value = "20:58:56"
$ws.Cells.Item(15, 3).Value = "20:00"
$ws.Cells.Item(15, 4).Value = "Bedroom"
$ws.Cells.Item(15, 5).Value = 50
$ws.Cells.Item(15, 6).Value = "Occupied"
$ws.Cells.Item(16, 1).Value = "'2026-02-01"
$ws.Cells.Item(16, 2).Value = "20:59:02"
$ws.Cells.Item(16, 3).Value = "20:00"
$ws.Cells.Item(16, 4).Value = "Bedroom"
$ws.Cells.Item(16, 5).Value = 58
$ws.Cells.Item(16, 6).Value = "Occupied"
$ws.Cells.Item(17, 1).Value = "'2026-02-01"
$ws.Cells.Item(17, 2).Value = "20:59:02"
$ws.Cells.Item(17, 3).Value = "20:00"
$ws.Cells.Item(17, 4).Value = "Bedroom"
$ws.Cells.Item(17, 5).Value = 52
$ws.Cells.Item(17, 6).Value = "Occupied"
$ws.Cells.Item(18, 1).Value = "'2026-02-01"
$ws.Cells.Item(18, 2).Value = "20:59:03"
$ws.Cells.Item(18, 3).Value = "20:00"
$ws.Cells.Item(18, 4).Value = "Bedroom"
$ws.Cells.Item(18, 5).Value = 50
$ws.Cells.Item(18, 6).Value = "Occupied"
$ws.Cells.Item(19, 1).Value = "'2026-02-01"
$ws.Cells.Item(19, 2).Value = "20:59:05"
$ws.Cells.Item(19, 3).Value = "20:00"
$ws.Cells.Item(19, 4).Value = "Bedroom"
$ws.Cells.Item(19, 5).Value = 68
$ws.Cells.Item(19, 6).Value = "Occupied"
$ws.Cells.Item(20, 1).Value = "'2026-02-01"
$ws.Cells.Item(20, 2).Value = "20:59:05"
$ws.Cells.Item(20, 3).Value = "20:00"
$ws.Cells.Item(20, 4).Value = "Bedroom"
$ws.Cells.Item(20, 5).Value = 50
$ws.Cells.Item(20, 6).Value = "Occupied"
$ws.Cells.Item(21, 1).Value = "'2026-02-01"
$ws.Cells.Item(21, 2).Value = "20:59:11"
$ws.Cells.Item(21, 3).Value = "20:00"
$ws.Cells.Item(21, 4).Value = "Bedroom"
$ws.Cells.Item(21, 5).Value = 49
$ws.Cells.Item(21, 6).Value = "Occupied"
$ws.Cells.Item(22, 1).Value = "'2026-02-01"
$ws.Cells.Item(22, 2).Value = "20:59:17"
$ws.Cells.Item(22, 3).Value = "20:00"
$ws.Cells.Item(22, 4).Value = "Bedroom"
$ws.Cells.Item(22, 5).Value = 50
$ws.Cells.Item(22, 6).Value = "Occupied"
$ws.Cells.Item(23, 1).Value = "'2026-02-01"
$ws.Cells.Item(23, 2).Value = "20:59:23"
$ws.Cells.Item(23, 3).Value = "20:00"
$ws.Cells.Item(23, 4).Value = "Bedroom"
$ws.Cells.Item(23, 5).Value = 120
$ws.Cells.Item(23, 6).Value = "Occupied"
$ws.Cells.Item(24, 1).Value = "'2026-02-01"
$ws.Cells.Item(24, 2).Value = "20:59:24"
$ws.Cells.Item(24, 3).Value = "20:00"
$ws.Cells.Item(24, 4).Value = "Bedroom"
$ws.Cells.Item(24, 5).Value = 139
$ws.Cells.Item(24, 6).Value = "Occupied"
$ws.Cells.Item(25, 1).Value = "'2026-02-01"
$ws.Cells.Item(25, 2).Value = "20:59:26"
$ws.Cells.Item(25, 3).Value = "20:00"
$ws.Cells.Item(25, 4).Value = "Bedroom"
$ws.Cells.Item(25, 5).Value = 50
$ws.Cells.Item(25, 6).Value = "Occupied"
$ws.Range("A13:A25").ClearFormats()
